$wb = $excel.ActiveWorkbook

$wsQuality = $wb.Worksheets.Item("quality_comparison")
$wsComp    = $wb.Worksheets.Item("computational_comparison")

# ------------------------------------------------------------------
# 1) New border-only cell styles for the (previously boxed) C1/D1
#    header cells: drop the left/right edges and keep only a
#    top+bottom rule (style "4"), or top+bottom+right (style "5").
#    Build these once on quality_comparison!C1 / D1 and then copy the
#    *formats only* onto the equivalent cells so every target cell
#    resolves to the same two new cellXfs entries instead of minting
#    a new one each time.
# ------------------------------------------------------------------
$c1 = $wsQuality.Range("C1")
$c1.ClearFormats()
$c1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$c1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

$d1 = $wsQuality.Range("D1")
$d1.ClearFormats()
$d1.Borders.Item(8).LineStyle = 1    # xlEdgeTop
$d1.Borders.Item(10).LineStyle = 1   # xlEdgeRight
$d1.Borders.Item(9).LineStyle = 1    # xlEdgeBottom

# Re-use the two new styles on computational_comparison's matching cells.
$c1.Copy()
$wsComp.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$wsComp.Range("F1").PasteSpecial(-4122)

$d1.Copy()
$wsComp.Range("D1").PasteSpecial(-4122)
$wsComp.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 2) Anonymize the "fedcore" header label -> "approach"
# ------------------------------------------------------------------
$wsQuality.Range("C2").Value = "approach"
$wsComp.Range("C2").Value = "approach"
$wsComp.Range("F2").Value = "approach"

# ------------------------------------------------------------------
# 3) Drop the stray empty inline-string cell G5 on computational_comparison
# ------------------------------------------------------------------
$wsComp.Range("G5").ClearContents()
